$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "GNG_TO-16502912156931021"
$wb.Worksheets.Item(2).Name = "NB_TO-16502912183219569"
$wb.Worksheets.Item(3).Name = "RS_TO-1650291218323956"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502912183869636"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502912184791565"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502912156400983.csv"
$ws1.Range("B3").Value = "GNG_stims-16502912156600986.csv"
$ws1.Range("B4").Value = "go_stims-1650291215661099.csv"
$ws1.Range("B5").Value = "GNG_stims-1650291215691101.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-1650291217113704.csv"
$ws2.Range("B3").Value = "OB-16502912167617052.csv"
$ws2.Range("B4").Value = "TB-16502912183049543.csv"
$ws2.Range("B5").Value = "ZB-match_7-16502912157631009.csv"
$ws2.Range("B6").Value = "OB-16502912163076992.csv"
$ws2.Range("B7").Value = "ZB-match_6-16502912161886969.csv"
$ws2.Range("B8").Value = "TB-16502912176929564.csv"
$ws2.Range("B9").Value = "ZB-match_2-16502912161547036.csv"
$ws2.Range("B10").Value = "TB-16502912181149566.csv"

# Sheet 3 (RS) - no cell content changes, only sheet name was updated above

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502912183379562.csv"
$ws4.Range("B3").Value = "ZM_stims-16502912183259568.csv"
$ws4.Range("B4").Value = "MM_stims-1650291218369957.csv"
$ws4.Range("B5").Value = "ZM_stims-16502912183389568.csv"
$ws4.Range("B6").Value = "MM_stims-16502912183859925.csv"
$ws4.Range("B7").Value = "ZM_stims-16502912183709593.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16502912184169555.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502912184331388.csv"
$ws5.Range("B4").Value = "SAT_stims-1650291218390995.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502912184631183.csv"
